# Updated Contributions with entities so far:
#  - Eumee:  add "Phase 5: Implemented Manage Game and Manage Console"
#  - Vanesa: add "Phase 5: Implemented Manage Inventory and Manage Product"
#  - Connar: add "Phase 5: Implemented Manage Store and Manage Employee"

$d = $word.ActiveDocument

function Get-ParagraphIndexContaining($rng) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Start -le $rng.Start -and $p.Range.End -ge $rng.End) {
            return $i
        }
    }
    return 0
}

# Inserts a brand-new "List Paragraph" bullet right after the paragraph that
# contains $searchText, populated with the run texts in $segments (each
# element becomes its own w:r, matching how Word naturally splits runs when
# text is typed/pasted incrementally).
function Insert-BulletAfter($searchText, $segments) {
    $rng = $d.Content
    $found = $rng.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Could not find anchor text: $searchText"
    }
    $idx = Get-ParagraphIndexContaining $rng
    if ($idx -eq 0) {
        throw "Could not resolve paragraph index for: $searchText"
    }
    $anchorPara = $d.Paragraphs.Item($idx)
    $anchorPara.Range.InsertParagraphAfter()

    $newPara = $d.Paragraphs.Item($idx + 1)
    $newPara.Range.Text = $segments[0]

    for ($j = 1; $j -lt $segments.Count; $j++) {
        $livePara = $d.Paragraphs.Item($idx + 1)
        $livePara.Range.InsertAfter($segments[$j])
    }
}

# 1) New bullet under "Eumee", after the "Coded the XAML..." bullet.
Insert-BulletAfter "Phase 5: Coded the XAML for the C# application" @("Phase 5: Implemented Manage ", "Game and Manage Console")

# 2) New bullet under "Vanesa", after the "Added necessary classes..." bullet.
Insert-BulletAfter "Phase 5: Added necessary classes and functions to C# application to connect to MySQL database and load data from entities into dataset" @("Phase 5: Implemented Manage Inventory and Manage Product")

# 3) New bullet under "Connar", after the last bullet ("...with XAML and code-behind").
Insert-BulletAfter "Phase 5: Created the load and save dataset functionality in the C# application with XAML and code-behind" @("Phase 5: Implemented Manage ", "Store", " and Manage ", "Employee")
